$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows whose Target cluster (column D) is "ECs" (deleted bottom-to-top)
$ws.Rows(8).Delete()
$ws.Rows(5).Delete()
$ws.Rows(2).Delete()

# Update remaining rows with recomputed TPM-derived values
# Row 2: ECs -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vtn"
$ws.Range("C2").Value = "Tnfrsf11b"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 12.056684
$ws.Range("H2").Value = 36.170052
$ws.Range("I2").Value = 0.06307822458376462
$ws.Range("J2").Value = 0.06307822458376462
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.183046666666667
$ws.Range("N2").Value = 3.54914
$ws.Range("O2").Value = 0.6222589862820888
$ws.Range("P2").Value = 0.6222589862820888
$ws.Range("Q2").Value = 14.26361981725333
$ws.Range("R2").Value = 128.37257835528
$ws.Range("S2").Value = 0.0392509920859673
$ws.Range("T2").Value = 0.0392509920859673

# Row 3: ECs -> MuSCs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vtn"
$ws.Range("C3").Value = "Tnfrsf11b"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 12.056684
$ws.Range("H3").Value = 36.170052
$ws.Range("I3").Value = 0.06307822458376462
$ws.Range("J3").Value = 0.06307822458376462
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.718166
$ws.Range("N3").Value = 2.154498
$ws.Range("O3").Value = 0.3777410137179113
$ws.Range("P3").Value = 0.3777410137179112
$ws.Range("Q3").Value = 8.658700521543999
$ws.Range("R3").Value = 77.928304693896
$ws.Range("S3").Value = 0.02382723249779732
$ws.Range("T3").Value = 0.02382723249779732

# Row 4: FAPs -> FAPs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Vtn"
$ws.Range("C4").Value = "Tnfrsf11b"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 25.140634
$ws.Range("H4").Value = 75.421902
$ws.Range("I4").Value = 0.1315309049843414
$ws.Range("J4").Value = 0.1315309049843414
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.183046666666667
$ws.Range("N4").Value = 3.54914
$ws.Range("O4").Value = 0.6222589862820888
$ws.Range("P4").Value = 0.6222589862820888
$ws.Range("Q4").Value = 29.74254325158667
$ws.Range("R4").Value = 267.68288926428
$ws.Range("S4").Value = 0.08184628760032199
$ws.Range("T4").Value = 0.08184628760032199

# Row 5: FAPs -> MuSCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Vtn"
$ws.Range("C5").Value = "Tnfrsf11b"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 25.140634
$ws.Range("H5").Value = 75.421902
$ws.Range("I5").Value = 0.1315309049843414
$ws.Range("J5").Value = 0.1315309049843414
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.718166
$ws.Range("N5").Value = 2.154498
$ws.Range("O5").Value = 0.3777410137179113
$ws.Range("P5").Value = 0.3777410137179112
$ws.Range("Q5").Value = 18.055148557244
$ws.Range("R5").Value = 162.496337015196
$ws.Range("S5").Value = 0.04968461738401938
$ws.Range("T5").Value = 0.04968461738401937

# Row 6: MuSCs -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Vtn"
$ws.Range("C6").Value = "Tnfrsf11b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 153.9412893333333
$ws.Range("H6").Value = 461.8238680000001
$ws.Range("I6").Value = 0.8053908704318941
$ws.Range("J6").Value = 0.8053908704318941
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.183046666666667
$ws.Range("N6").Value = 3.54914
$ws.Range("O6").Value = 0.6222589862820888
$ws.Range("P6").Value = 0.6222589862820888
$ws.Range("Q6").Value = 182.1197292081689
$ws.Range("R6").Value = 1639.07756287352
$ws.Range("S6").Value = 0.5011617065957995
$ws.Range("T6").Value = 0.5011617065957995

# Row 7: MuSCs -> MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Vtn"
$ws.Range("C7").Value = "Tnfrsf11b"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 153.9412893333333
$ws.Range("H7").Value = 461.8238680000001
$ws.Range("I7").Value = 0.8053908704318941
$ws.Range("J7").Value = 0.8053908704318941
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.718166
$ws.Range("N7").Value = 2.154498
$ws.Range("O7").Value = 0.3777410137179113
$ws.Range("P7").Value = 0.3777410137179112
$ws.Range("Q7").Value = 110.5553999953627
$ws.Range("R7").Value = 994.9985999582641
$ws.Range("S7").Value = 0.3042291638360946
$ws.Range("T7").Value = 0.3042291638360946

